$d = $word.ActiveDocument
$p = $d.Paragraphs(2)
$p.Style = "List Paragraph"
